$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Price (D) and Volume(1h) (E) columns with the latest crypto snapshot values.
# Column D holds price text (e.g. "66.578.65") that must stay as text, so we
# force a Text number format before assigning the value to avoid Excel
# auto-converting numeric-looking strings (like "609.52") into real numbers.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '66.578.65'
$ws.Cells.Item(2, 5).Value = '  +0.69%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.600.66'
$ws.Cells.Item(3, 5).Value = '  +1.34%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '609.52'
$ws.Cells.Item(5, 5).Value = '  +0.57%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '149.05'
$ws.Cells.Item(6, 5).Value = '  +3.10%  '
$ws.Cells.Item(7, 5).Value = '  +0.19%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.488'
$ws.Cells.Item(8, 5).Value = '  -0.74%  '
$ws.Cells.Item(9, 5).Value = '  +1.91%  '
$ws.Cells.Item(10, 5).Value = '  -0.13%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.416'
$ws.Cells.Item(11, 5).Value = '  +0.67%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '4.214.22'
$ws.Cells.Item(12, 5).Value = '  +1.43%  '
$ws.Cells.Item(13, 5).Value = '  +1.10%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '29.78'
$ws.Cells.Item(14, 5).Value = '  -0.96%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.593.83'
$ws.Cells.Item(15, 5).Value = '  +1.04%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '66.659.46'
$ws.Cells.Item(16, 5).Value = '  +0.75%  '
$ws.Cells.Item(17, 5).Value = '  +0.74%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '11.55'
$ws.Cells.Item(18, 5).Value = '  +1.92%  '
$ws.Cells.Item(19, 5).Value = '  +3.25%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '15.12'
$ws.Cells.Item(20, 5).Value = '  +1.79%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '428.37'
$ws.Cells.Item(21, 5).Value = '  -0.46%  '
$ws.Cells.Item(22, 5).Value = '  +1.47%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '78.94'
$ws.Cells.Item(23, 5).Value = '  -0.31%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.751.28'
$ws.Cells.Item(24, 5).Value = '  +1.59%  '
$ws.Cells.Item(25, 5).Value = '  -0.03%  '
$ws.Cells.Item(26, 5).Value = '  +4.20%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '8.31'
$ws.Cells.Item(27, 5).Value = '  +4.32%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.46'
$ws.Cells.Item(28, 5).Value = '  +4.13%  '
$ws.Cells.Item(29, 5).Value = '  +0.29%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  -0.10%  '
$ws.Cells.Item(31, 5).Value = '  +0.83%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.598.82'
$ws.Cells.Item(32, 5).Value = '  +1.48%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.158'
$ws.Cells.Item(33, 5).Value = '  +3.95%  '
$ws.Cells.Item(34, 5).Value = '  -0.13%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '7.89'
$ws.Cells.Item(35, 5).Value = '  +0.11%  '
$ws.Cells.Item(36, 5).Value = '  -0.02%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.67'
$ws.Cells.Item(37, 5).Value = '  +1.03%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.70'
$ws.Cells.Item(38, 5).Value = '  -2.33%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '178.04'
$ws.Cells.Item(39, 5).Value = '  +2.38%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0857'
$ws.Cells.Item(40, 5).Value = '  +0.75%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '5.25'
$ws.Cells.Item(41, 5).Value = '  +0.79%  '
$ws.Cells.Item(42, 5).Value = '  +0.83%  '
$ws.Cells.Item(43, 5).Value = '  -0.86%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.56'
$ws.Cells.Item(44, 5).Value = '  +9.83%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.00'
$ws.Cells.Item(45, 5).Value = '  +0.16%  '
$ws.Cells.Item(46, 5).Value = '  -2.10%  '
$ws.Cells.Item(47, 5).Value = '  -1.71%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '24.10'
$ws.Cells.Item(48, 5).Value = '  +2.04%  '
$ws.Cells.Item(49, 5).Value = '  +1.10%  '
$ws.Cells.Item(50, 5).Value = '  +1.48%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.432.16'
$ws.Cells.Item(51, 5).Value = '  +5.70%  '
